$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 52 (pushes the old row 52.."idLevelTitle000" block etc. down by one,
# i.e. old rows 52-121 become 53-122) for the new "idShowAd" localization entry.
$ws.Rows(52).Insert()

# Populate the three new cells with the localization id, English copy and Russian copy.
$ws.Range("A52").Value = "idShowAd"
$ws.Range("B52").Value = "This button lets you view an ad and earn some RotoCoins!"
$ws.Range("C52").Value = "Нажав эту кнопку, Вы посмотрите рекламу и получите немного РотоКойнов!"

# Match the taller row height used for this entry (wrapped two-line text).
$ws.Rows(52).RowHeight = 26.4

# Update the view state to reflect where the author ended up after the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 42
$win.ScrollColumn = 1
$ws.Range("A52").Select() | Out-Null
